$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank columns in front of the existing A:B data so it shifts to
# E:F. Using EntireColumn.Insert() (rather than clearing + rewriting values)
# carries the cell contents, types, and the column B custom width along to
# their new home (E:F) exactly, byte-for-byte, instead of recomputing them.
$ws.Range("A1:D1").EntireColumn.Insert()

# Match the saved selection state (whole-column selection anchored at A1).
$ws.Range("A1:A1048576").Select()
